$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.966.36"
$ws.Range("E2").Value = "  +1.67%  "

$ws.Range("D3").Value = "1.648.59"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'213.53"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D8").Value = "'23.61"
$ws.Range("E8").Value = "  +3.91%  "

$ws.Range("E9").Value = "  +1.37%  "

$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").Value = "1.881.58"
$ws.Range("E12").Value = "  +1.98%  "

$ws.Range("D13").Value = "1.654.76"
$ws.Range("E13").Value = "  +2.32%  "

$ws.Range("D14").Value = "'4.09"
$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("D16").Value = "'65.73"
$ws.Range("E16").Value = "  +1.38%  "

$ws.Range("D17").Value = "27.957.06"
$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("D18").Value = "'232.09"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").Value = "'7.68"
$ws.Range("E19").Value = "  +2.33%  "

$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "'10.71"
$ws.Range("E22").Value = "  +6.22%  "

$ws.Range("E23").Value = "  +2.62%  "

$ws.Range("E24").Value = "  +3.05%  "

$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("D26").Value = "'6.93"
$ws.Range("E26").Value = "  +1.88%  "

$ws.Range("D27").Value = "'15.74"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "'1.19"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("D32").Value = "'3.34"
$ws.Range("E32").Value = "  +2.47%  "

$ws.Range("D33").Value = "1.454.56"
$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("E34").Value = "  +2.33%  "

$ws.Range("E35").Value = "  +2.10%  "

$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("E37").Value = "  +3.64%  "

$ws.Range("E38").Value = "  +1.16%  "

$ws.Range("D39").Value = "'0.563"
$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("D40").Value = "'0.916"
$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("D41").Value = "'69.46"
$ws.Range("E41").Value = "  +0.59%  "

$ws.Range("E42").Value = "  +1.85%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("E45").Value = "  +1.63%  "

$ws.Range("D46").Value = "'5.39"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("E47").Value = "  +6.11%  "

$ws.Range("D48").Value = "1.791.36"
$ws.Range("E48").Value = "  +1.88%  "

$ws.Range("D49").Value = "'89.04"
$ws.Range("E49").Value = "  +3.04%  "

$ws.Range("E51").Value = "  +1.53%  "
